$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.245.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.71%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("E7").Value = "  -0.70%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3742"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07163"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8944"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.881.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.314"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9989"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008518"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.264.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.014"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.106.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.73%  "

$ws.Range("E23").Value = "  -3.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.483"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.94%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.835"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.74%  "

$ws.Range("E27").Value = "  -1.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.094"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.671"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.692"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09267"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05139"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.83%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.087"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.164"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7292"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.139"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02035"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.527"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5326"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.544"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.366"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.59%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1478"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4643"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9992"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.567"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.28%  "
